$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A35:C35").ClearContents()
$ws.Range("D35:T35").ClearContents()

$ws.Range("Y24").Select()
